# 自动更新价格数据: insert a new row for 2026-01-04 at the top of the
# data (row 2), pushing all existing date rows down by one. The new
# row carries the same commodity values as the most recent row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts rows 2..45 -> 3..46)
$ws.Rows.Item(2).Insert()

# Excel's row insert copies formatting down from the row above (the bold
# header row); the data rows in this sheet carry no explicit style, so
# strip whatever got copied in.
$ws.Range("A2:D2").ClearFormats()

# Populate the newly inserted row with the latest date + carried-over values.
# The date column stores plain text (not a real date) elsewhere in the
# sheet, so lead with an apostrophe to force text entry and keep Excel
# from auto-converting it to a date serial number.
$ws.Cells.Item(2, 1).Value = "'2026-01-04"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# The quote-prefix entry stamps a "quotePrefix" style on the cell; clear
# formats again so the new row matches the unstyled look of the other
# data rows.
$ws.Range("A2:D2").ClearFormats()
